# Applies the "update asset CSV data and modify dummy Excel data" edit:
#  - Re-points 9 Asset rows' sub_site_id (column P) from 25 -> 27
#  - Switches the active/selected sheet from "User" to "Asset" and updates
#    the Asset sheet's selection to P38

$wb = $excel.ActiveWorkbook

# --- 1. Update Asset sheet sub_site_id values (P column) from 25 to 27 ---
$assetSheet = $wb.Worksheets.Item("Asset")
$rowsToUpdate = @(2, 8, 9, 10, 11, 12, 25, 26, 27)
foreach ($r in $rowsToUpdate) {
    $assetSheet.Cells.Item($r, 16).Value2 = 27
}

# --- 2. Move the active tab from "User" to "Asset" and update selection ---
$assetSheet.Activate()
$assetSheet.Range("P38").Select()
